# Fruta / hortaliza, semanal
# New weekly price record inserted as row 47 ("Terminal Hortofrutícola Agro
# Chillán" - Mango, fecha 2022-02-11), pushing the previously-existing
# rows 47-72 down to 48-73 (dimension grows from A1:T72 to A1:T73).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 47; Excel shifts rows 47:72 -> 48:73
# and carries the row-47 formatting (date style on column D) down with it.
$ws.Rows.Item(47).Insert()

$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C47").Value = "Ñuble"
$ws.Range("D47").Value = (Get-Date -Year 2022 -Month 2 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E47").Value = 16
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100108
$ws.Range("H47").Value = "Tropicales y subtropicales"
$ws.Range("I47").Value = 100108002
$ws.Range("J47").Value = "Mango"
$ws.Range("K47").Value = "Sin especificar"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 60
$ws.Range("N47").Value = 7000
$ws.Range("O47").Value = 7500
$ws.Range("P47").Value = 7250
$ws.Range("Q47").Value = "$/bandeja 4 kilos"
$ws.Range("R47").Value = "Perú"
$ws.Range("S47").Value = 1812
$ws.Range("T47").Value = 4

Write-Host "Inserted new row 47; sheet now spans" $ws.UsedRange.Rows.Count "rows"
